$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.859.75"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "2.252.97"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'307.66"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'96.27"
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("D8").Value = "'1.01"
$ws.Range("D10").Value = "'35.53"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'7.28"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "2.596.07"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "2.311.13"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "'0.844"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "'13.63"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "44.562.11"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.33"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'12.06"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'65.65"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'239.02"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'37.39"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'152.78"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "'3.12"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("D39").Value = "'14.98"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'3.84"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "1.844.16"
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("E45").Value = "  +17.86%  "
$ws.Range("D46").Value = "'0.193"
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("D47").Value = "'80.09"
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("D48").Value = "'99.72"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'70.67"
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'55.21"
$ws.Range("E51").Value = "  +2.67%  "
